# "layout added" - drop the x/y scatter-plot coordinate columns from the
# Sheet2 "nodes" table; the trailing `score` column slides left to take
# their place (F), and the now-unused "x"/"y" shared strings disappear.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Columns F ("x") and G ("y") go away entirely; H ("score") shifts left
# into F. EntireColumn.Delete shifts everything after G left by two and
# drops the dimension/shared-strings bookkeeping for us.
$ws.Range("F1:G1").EntireColumn.Delete()

# Selection now covers the (new) last two columns, matching the
# post-edit view state.
$ws.Range("F1:G1048576").Select()
